# Fill in the second data row (row 3) of the UserData sheet with a new
# signup test-data record, mirroring the pattern already used in row 2.
# The data is entered in the same order a tester would tab through the
# signup form so that newly-introduced shared strings line up with the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserData")
$ws.Activate() | Out-Null

# id / title
$ws.Range("A3").Value = "standard"
$ws.Range("B3").Value = "Mrs"

# name (lastName typed right after title, then firstName)
$ws.Range("D3").Value = "a"
$ws.Range("C3").Value = "A"

# email + password
$ws.Range("E3").Value = "hello@yahoo.com"
$ws.Range("F3").Value = "qew2"

# date of birth
$ws.Range("G3").Value = "17"
$ws.Range("H3").Value = "2"
$ws.Range("I3").Value = "1995"

# address block
$ws.Range("J3").Value = "34524dsada"
$ws.Range("K3").Value = "Cairo"
$ws.Range("L3").Value = "Nasr City"
$ws.Range("M3").Value = "54353"

# mobile number, then shouldSucceed flag, then country
$ws.Range("O3").Value = "01232434243"

# "TRUE" must land in the cell as literal text (matching the column's
# Text format) instead of being auto-coerced into a native boolean, so
# build it as a formula result and paste the computed value back in.
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Formula = '="T"&"RUE"'
$ws.Range("P3").Copy() | Out-Null
$ws.Range("P3").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("N3").Value = "India"

# Turn the new e-mail address into a mailto hyperlink, like row 2's E2,
# then restore the Hyperlink cell style/text format that Hyperlinks.Add
# resets so the cell keeps using the original shared "Hyperlink" style.
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:hello@yahoo.com") | Out-Null
$ws.Range("E3").NumberFormat = "@"

# Reflect the cursor/scroll position that was left in the sheet after
# entering the new row: scrolled right so column D is left-most, with
# the new country cell (N3) selected.
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("N3").Select() | Out-Null
